$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (nocsFrameCounts -> FrameCounts)
$ws.Name = "FrameCounts"

# 2. Section header text update (shared string re-used by A5): Circuits -> Level 1
$ws.Cells.Item(5, 1).Value = "Level 1"

# 3. Header row labels (row 1): V6/V5 -> Mine/Andymac
$ws.Cells.Item(1, 2).Value = "Mine"
$ws.Cells.Item(1, 3).Value = "Andymac"

# 4. Replace the data table (rows 6-12) with Sonikkustars' new 4-frame-faster numbers
$ws.Cells.Item(6, 1).Value = "Batman appears"
$ws.Cells.Item(6, 2).Value = 450
$ws.Cells.Item(6, 3).Value = 450

$ws.Cells.Item(7, 1).Value = "X = 210"
$ws.Cells.Item(7, 2).Value = 570
$ws.Cells.Item(7, 3).Value = 622

$ws.Cells.Item(8, 1).Value = "X = 579"
$ws.Cells.Item(8, 2).Value = 690
$ws.Cells.Item(8, 3).Value = 741

$ws.Cells.Item(9, 1).Value = "X = 901"
$ws.Cells.Item(9, 2).Value = 844
$ws.Cells.Item(9, 3).Value = 900

$ws.Cells.Item(10, 1).Value = "X = 1129"
$ws.Cells.Item(10, 2).Value = 988
$ws.Cells.Item(10, 3).Value = 1040

$ws.Cells.Item(11, 1).Value = "X = 1252"
$ws.Cells.Item(11, 2).Value = 1281
$ws.Cells.Item(11, 3).Value = 1330

$ws.Cells.Item(12, 1).Value = "Screen 2"
$ws.Cells.Item(12, 2).Value = 1666
$ws.Cells.Item(12, 3).Value = 1677

# 5. The remaining old rows (13-32) lose all their old place/frame data - wipe them
#    out completely (not just clear contents) so the empty cells disappear from
#    the sheet, matching the cleaned-up worksheet.
$ws.Range("A13:C32").Clear()

# 6. Column C now needs a touch more width to fit the new numbers.
$ws.Columns.Item(3).ColumnWidth = 9.5

# 7. Refresh the frozen-pane view: scrolled back to the top, selection sitting on B13.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B13").Select()
